# Update the "Förändrad" (Changed) date column (C) for rows 2-26
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
